$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the activity text in A5: "market for electricity, low voltage"
# -> "market group for electricity, low voltage"
$ws.Cells.Item(5, 1).Value = "market group for electricity, low voltage"

# Updated sensitivity results in column B (values are stored as text, not
# numbers, in the source workbook - force Text format before assignment so
# Excel does not auto-convert the numeric-looking strings into real numbers).
$updates = @{
    2  = "0.6003256557974975"
    3  = "0.7724279569836534"
    4  = "1.5930136661490186"
    5  = "15.910408342898313"
    6  = "15.48664901692736"
    7  = "1.984545550413917"
    8  = "4.036532437006933"
    9  = "0.2766790808148134"
    10 = "4.938014071436373"
    11 = "1.2369982059045173"
    12 = "0.2861877148894056"
    13 = "0.002426815732611057"
    14 = "0.01538042279581438"
    15 = "2.47225783771006"
    16 = "0.002970614600214056"
    17 = "0.08667283050730451"
    18 = "-0.14273056182503313"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    # Force Text storage so the numeric-looking string is kept as a literal
    # string (matching the source file) instead of being auto-converted to
    # a real number, then restore the default "Normal" style so no new
    # number-format style gets left attached to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}
